$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = 9
$ws.Range("F3").Value = -5
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = 6
$ws.Range("F13").Value = -2
$ws.Range("F14").Value = -4
